$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A3").Value = "Vrinda"
$ws.Range("B3").Value = 97

$ws.Range("B3").Select()
